$wb = $excel.ActiveWorkbook

# --- Sheet 1: "Pediatric VFC Vaccine " ---
$ws1 = $wb.Worksheets.Item(1)

# "DTaP" -> "DTaP/"
$ws1.Range("A2").Value = "DTaP/"
$ws1.Range("A3").Value = "DTaP/"
$ws1.Range("A4").Value = "DTaP/"
$ws1.Range("A5").Value = "DTaP/"

# "DTaP " -> "DTaP-Hib "
$ws1.Range("A8").Value = "DTaP-Hib "

# "Hepatitis B^" -> "Hepatitis B-Hib"
$ws1.Range("A10").Value = "Hepatitis B-Hib"

# "Hepatitis A-Hepatitis B 18 only^" -> "Hepatitis A-Hepatitis B 18 only"
$ws1.Range("A15").Value = "Hepatitis A-Hepatitis B 18 only"
$ws1.Range("A16").Value = "Hepatitis A-Hepatitis B 18 only"

# "Hepatitis B PediatricAdolescent" -> "Hepatitis B Pediatric/Adolescent"
$ws1.Range("A17").Value = "Hepatitis B Pediatric/Adolescent"
$ws1.Range("A18").Value = "Hepatitis B Pediatric/Adolescent"
$ws1.Range("A19").Value = "Hepatitis B Pediatric/Adolescent"
$ws1.Range("A20").Value = "Hepatitis B Pediatric/Adolescent"

# "MMR" -> "MMR/"
$ws1.Range("A27").Value = "MMR/"

# "Pneumococcal7-valent (Pediatric)" -> "Pneumococcal 7-valent (Pediatric)"
$ws1.Range("A28").Value = "Pneumococcal 7-valent (Pediatric)"

# "Tetanus  Diphtheria Toxoids^" -> "Tetanus  Diphtheria Toxoids"
$ws1.Range("A29").Value = "Tetanus  Diphtheria Toxoids"

# "Tetanus Toxoid, Reduced Diphtheria Toxoid and Acellular Pertussis" -> "...Pertussis/"
$ws1.Range("A30").Value = "Tetanus Toxoid, Reduced Diphtheria Toxoid and Acellular Pertussis/"
$ws1.Range("A31").Value = "Tetanus Toxoid, Reduced Diphtheria Toxoid and Acellular Pertussis/"
$ws1.Range("A32").Value = "Tetanus Toxoid, Reduced Diphtheria Toxoid and Acellular Pertussis/"

# Split the combined packaging text into two separate rows/strings:
# D30: "10 pack - 1 dose vials 5 pack - 1 dose TL syringes, No Needle " -> "10 pack - 1 dose vials "
# D31: same original text -> new string "5 pack - 1 dose TL syringes, No Needle "
$ws1.Range("D30").Value = "10 pack - 1 dose vials "
$ws1.Range("D31").Value = "5 pack - 1 dose TL syringes, No Needle "

# --- Sheet 2: "Adult Vaccine " ---
$ws2 = $wb.Worksheets.Item(2)

# "Hepatitis A-Hepatitis B Adult^" -> "Hepatitis A-Hepatitis B Adult"
$ws2.Range("A6").Value = "Hepatitis A-Hepatitis B Adult"
$ws2.Range("A7").Value = "Hepatitis A-Hepatitis B Adult"

# "Tetanus  Diphtheria Toxoids^" -> "Tetanus  Diphtheria Toxoids"
$ws2.Range("A13").Value = "Tetanus  Diphtheria Toxoids"
